$d = $word.ActiveDocument

$d.Content.Find.Execute("Deliverables:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Potential Deliverables:", 2)

$d.Save()
